$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '27.707.98'
Set-TextValue 'E2' '  -0.83%  '
Set-TextValue 'D3' '1.585.35'
Set-TextValue 'E3' '  -3.05%  '
Set-TextValue 'E4' '  +0.34%  '
Set-TextValue 'D5' '206.86'
Set-TextValue 'E5' '  -2.23%  '
Set-TextValue 'D6' '0.508'
Set-TextValue 'E6' '  -2.28%  '
Set-TextValue 'E7' '  +0.36%  '
Set-TextValue 'D8' '22.29'
Set-TextValue 'E8' '  -4.86%  '
Set-TextValue 'E9' '  -1.59%  '
Set-TextValue 'D10' '0.0592'
Set-TextValue 'E10' '  -3.04%  '
Set-TextValue 'D11' '0.0868'
Set-TextValue 'E11' '  -1.56%  '
Set-TextValue 'D12' '1.808.26'
Set-TextValue 'E12' '  -3.18%  '
Set-TextValue 'D13' '1.569.47'
Set-TextValue 'E13' '  -4.40%  '
Set-TextValue 'E14' '  -3.70%  '
Set-TextValue 'D15' '0.531'
Set-TextValue 'E15' '  -5.86%  '
Set-TextValue 'D16' '27.662.77'
Set-TextValue 'E16' '  -1.01%  '
Set-TextValue 'D17' '63.30'
Set-TextValue 'E17' '  -2.98%  '
Set-TextValue 'D18' '220.06'
Set-TextValue 'E18' '  -4.18%  '
Set-TextValue 'D19' '0.0₃0694'
Set-TextValue 'E19' '  -3.62%  '
Set-TextValue 'E20' '  -6.51%  '
Set-TextValue 'E21' '  +0.39%  '
Set-TextValue 'E22' '  -4.97%  '
Set-TextValue 'D23' '9.51'
Set-TextValue 'E23' '  -6.19%  '
Set-TextValue 'E24' '  -4.95%  '
Set-TextValue 'D25' '154.18'
Set-TextValue 'E25' '  -1.26%  '
Set-TextValue 'D26' '6.79'
Set-TextValue 'E26' '  -2.75%  '
Set-TextValue 'E27' '  +0.37%  '
Set-TextValue 'D28' '15.16'
Set-TextValue 'E28' '  -2.53%  '
Set-TextValue 'E29' '  -3.94%  '
Set-TextValue 'D30' '1.15'
Set-TextValue 'E30' '  -2.32%  '
Set-TextValue 'D31' '0.0465'
Set-TextValue 'E31' '  -3.46%  '
Set-TextValue 'E32' '  -5.30%  '
Set-TextValue 'D33' '1.387.60'
Set-TextValue 'E33' '  -0.98%  '
Set-TextValue 'E34' '  -5.50%  '
Set-TextValue 'E35' '  -5.39%  '
Set-TextValue 'D36' '0.969'
Set-TextValue 'E36' '  -4.47%  '
Set-TextValue 'E38' '  -3.43%  '
Set-TextValue 'E39' '  -3.70%  '
Set-TextValue 'D40' '0.821'
Set-TextValue 'E40' '  -3.70%  '
Set-TextValue 'E41' '  +0.39%  '
Set-TextValue 'D42' '0.977'
Set-TextValue 'E42' '  -3.66%  '
Set-TextValue 'B43' 'Aave'
Set-TextValue 'C43' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D43' '63.63'
Set-TextValue 'E43' '  -3.78%  '
Set-TextValue 'B44' 'MXToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D44' '2.17'
Set-TextValue 'E44' '  +1.38%  '
Set-TextValue 'E45' '  -4.43%  '
Set-TextValue 'D46' '5.23'
Set-TextValue 'E46' '  -4.21%  '
Set-TextValue 'D47' '1.719.39'
Set-TextValue 'D48' '88.11'
Set-TextValue 'E48' '  -0.48%  '
Set-TextValue 'E49' '  -2.20%  '
Set-TextValue 'D50' '0.0975'
Set-TextValue 'E50' '  -4.69%  '
Set-TextValue 'D51' '0.0498'
Set-TextValue 'E51' '  -1.44%  '
